$d = $word.ActiveDocument

# Locate the paragraph containing the "Ver no Jupiter..." text and the one
# before it (blank) plus the "(c) 2020..." paragraph after it, and delete
# the whole range (including paragraph marks) so the three paragraphs are
# removed entirely.

$start = $null
$end = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        # the blank paragraph right before this one
        $prev = $d.Paragraphs.Item($i - 1)
        $start = $prev.Range.Start
        # the copyright paragraph right after this one
        $next = $d.Paragraphs.Item($i + 1)
        $end = $next.Range.End
        break
    }
}

if ($start -ne $null -and $end -ne $null) {
    $r = $d.Range($start, $end)
    $r.Delete()
}
